$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy the style of an existing header cell (AC1)
# so they get the same bold/bordered/centered formatting, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-49: team record values (Wins=95, Losses=67, Ties=0)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 95   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 67   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
